$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GenericDataset")

# Duplicate the GenericDataset sheet (mirrors Excel's "Move or Copy > Create a copy")
$ws.Copy($null, $ws)
$newWs = $wb.Worksheets.Item($ws.Index + 1)
$newWs.Name = "FieldNotes"

# Update the reference_type values on the new sheet
$newWs.Range("A1").Value = "FieldNotes"
$newWs.Range("A2").Value = "FieldNotes"

# Rebuild the validations in list(C1), date(E1:F1), list(A1:A2) order -
# matches how Excel re-serializes dataValidations after the sheet copy
$newWs.Range("C1").Validation.Delete()
$newWs.Range("C1").Validation.Add(3, 1, 1, """Yes, No""")
$newWs.Range("C1").Validation.InCellDropdown = $true
$newWs.Range("C1").Validation.ShowInput = $true
$newWs.Range("C1").Validation.ShowError = $true

$newWs.Range("E1:F1").Validation.Delete()
$newWs.Range("E1:F1").Validation.Add(4, 5, 1, "18264")
$newWs.Range("E1:F1").Validation.ShowInput = $true
$newWs.Range("E1:F1").Validation.ShowError = $true

# Update the data validation list on the new sheet for A1:A2 to restrict to FieldNotes
$newWs.Range("A1:A2").Validation.Delete()
$newWs.Range("A1:A2").Validation.Add(3, 1, 1, """FieldNotes""")
$newWs.Range("A1:A2").Validation.InCellDropdown = $true
$newWs.Range("A1:A2").Validation.ShowInput = $true
$newWs.Range("A1:A2").Validation.ShowError = $true

# Restore selection/active-cell states to mirror a typical copy operation
$ws.Activate()
$ws.Rows("1:2").Select()

$newWs.Activate()
$newWs.Range("K2").Select()
